$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37771
$ws.Range("D2").Value = 54625680
$ws.Range("C3").Value = 91022
$ws.Range("D3").Value = 133426584
$ws.Range("C4").Value = 31190
$ws.Range("D4").Value = 46191349
$ws.Range("C5").Value = 8702
$ws.Range("D5").Value = 12933563
$ws.Range("C6").Value = 1998
$ws.Range("D6").Value = 2968971
$ws.Range("C12").Value = 41362
$ws.Range("D12").Value = 56116411
$ws.Range("C13").Value = 9659
$ws.Range("D13").Value = 13970443
$ws.Range("C14").Value = 25954
$ws.Range("D14").Value = 38062884
$ws.Range("C15").Value = 8310
$ws.Range("D15").Value = 12332824
$ws.Range("C17").Value = 417
$ws.Range("D17").Value = 614623
$ws.Range("C20").Value = 10230
$ws.Range("D20").Value = 13543346
$ws.Range("C21").Value = 13388
$ws.Range("D21").Value = 19331137
$ws.Range("C22").Value = 31660
$ws.Range("D22").Value = 46460395
$ws.Range("C23").Value = 10223
$ws.Range("D23").Value = 15196491
$ws.Range("C27").Value = 11693
$ws.Range("D27").Value = 15618259
$ws.Range("C28").Value = 7643
$ws.Range("D28").Value = 11071117
$ws.Range("C29").Value = 22486
$ws.Range("D29").Value = 33006539
$ws.Range("C31").Value = 1960
$ws.Range("D31").Value = 2924499
$ws.Range("C34").Value = 8313
$ws.Range("D34").Value = 10979839
$ws.Range("C35").Value = 3247
$ws.Range("D35").Value = 4687194
$ws.Range("C36").Value = 7836
$ws.Range("D36").Value = 11443934
$ws.Range("C41").Value = 2473
$ws.Range("D41").Value = 3342853
$ws.Range("C42").Value = 17247
$ws.Range("D42").Value = 24937878
$ws.Range("C43").Value = 51132
$ws.Range("D43").Value = 74959865
$ws.Range("C44").Value = 19022
$ws.Range("D44").Value = 28255443
$ws.Range("C45").Value = 5610
$ws.Range("D45").Value = 8354177
$ws.Range("C50").Value = 16707
$ws.Range("D50").Value = 22237124
$ws.Range("C52").Value = 6906
$ws.Range("D52").Value = 10152079
$ws.Range("C53").Value = 2350
$ws.Range("D53").Value = 3509918
$ws.Range("C57").Value = 6993
$ws.Range("D57").Value = 9611817
$ws.Range("C58").Value = 951
$ws.Range("D58").Value = 1396079
$ws.Range("C59").Value = 2386
$ws.Range("D59").Value = 3537837
$ws.Range("C60").Value = 948
$ws.Range("D60").Value = 1411501
$ws.Range("C61").Value = 323
$ws.Range("D61").Value = 484258
$ws.Range("C63").Value = 18
$ws.Range("D63").Value = 27000
$ws.Range("C64").Value = 1400
$ws.Range("D64").Value = 1970306
$ws.Range("C65").Value = 15366
$ws.Range("D65").Value = 22195306
$ws.Range("C66").Value = 44687
$ws.Range("D66").Value = 65393525
$ws.Range("C67").Value = 15704
$ws.Range("D67").Value = 23338101
$ws.Range("C68").Value = 4569
$ws.Range("D68").Value = 6804703
$ws.Range("C69").Value = 925
$ws.Range("D69").Value = 1375668
$ws.Range("C73").Value = 15091
$ws.Range("D73").Value = 19895819
$ws.Range("C74").Value = 51417
$ws.Range("D74").Value = 74819478
$ws.Range("C75").Value = 146113
$ws.Range("D75").Value = 215255534
$ws.Range("C76").Value = 63642
$ws.Range("D76").Value = 94835657
$ws.Range("C77").Value = 20348
$ws.Range("D77").Value = 30402331
$ws.Range("C78").Value = 4820
$ws.Range("D78").Value = 7199043
$ws.Range("C85").Value = 50849
$ws.Range("D85").Value = 69166028
$ws.Range("C86").Value = 4606
$ws.Range("D86").Value = 6673050
$ws.Range("C87").Value = 11569
$ws.Range("D87").Value = 16996319
$ws.Range("C88").Value = 3884
$ws.Range("D88").Value = 5788583
$ws.Range("C93").Value = 5416
$ws.Range("D93").Value = 7281775
$ws.Range("C95").Value = 5166
$ws.Range("D95").Value = 7607743
$ws.Range("C96").Value = 1940
$ws.Range("D96").Value = 2889937
$ws.Range("C98").Value = 185
$ws.Range("D98").Value = 276613
$ws.Range("C99").Value = 19
$ws.Range("D99").Value = 28500
$ws.Range("C101").Value = 3566
$ws.Range("D101").Value = 4719108
$ws.Range("C102").Value = 604
$ws.Range("D102").Value = 899664
$ws.Range("C104").Value = 130
$ws.Range("D104").Value = 195000
$ws.Range("C107").Value = 10759
$ws.Range("D107").Value = 15607398
$ws.Range("C108").Value = 29208
$ws.Range("D108").Value = 42911810
$ws.Range("C109").Value = 9779
$ws.Range("D109").Value = 14541412
$ws.Range("C110").Value = 2686
$ws.Range("D110").Value = 4005207
$ws.Range("C111").Value = 491
$ws.Range("D111").Value = 731546
$ws.Range("C114").Value = 9794
$ws.Range("D114").Value = 12937606
$ws.Range("C115").Value = 30461
$ws.Range("D115").Value = 43924044
$ws.Range("C116").Value = 66157
$ws.Range("D116").Value = 96815230
$ws.Range("C117").Value = 21367
$ws.Range("D117").Value = 31754940
$ws.Range("C119").Value = 1121
$ws.Range("D119").Value = 1675271
$ws.Range("C121").Value = 13
$ws.Range("D121").Value = 19500
$ws.Range("C124").Value = 25851
$ws.Range("D124").Value = 34527269
$ws.Range("C125").Value = 35989
$ws.Range("D125").Value = 51940434
$ws.Range("C126").Value = 76786
$ws.Range("D126").Value = 112281580
$ws.Range("C127").Value = 23849
$ws.Range("D127").Value = 35394909
$ws.Range("C128").Value = 6393
$ws.Range("D128").Value = 9500238
$ws.Range("C129").Value = 1237
$ws.Range("D129").Value = 1839911
$ws.Range("C133").Value = 31822
$ws.Range("D133").Value = 42253084
$ws.Range("C134").Value = 13220
$ws.Range("D134").Value = 19135743
$ws.Range("C135").Value = 32321
$ws.Range("D135").Value = 47472152
$ws.Range("C136").Value = 11470
$ws.Range("D136").Value = 17042292
$ws.Range("C141").Value = 10817
$ws.Range("D141").Value = 14425134
$ws.Range("C142").Value = 35050
$ws.Range("D142").Value = 50614490
$ws.Range("C143").Value = 81195
$ws.Range("D143").Value = 118960416
$ws.Range("C144").Value = 24333
$ws.Range("D144").Value = 36153275
$ws.Range("C145").Value = 6388
$ws.Range("D145").Value = 9531567
$ws.Range("C146").Value = 1435
$ws.Range("D146").Value = 2134730
$ws.Range("C149").Value = 29201
$ws.Range("D149").Value = 39388986
